$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "Contact No"
$ws.Range("D1").Value = "Solution Required"
$ws.Range("E1").Value = "Message"
